$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'29.505.53"
$c.Style = "Normal"
$c = $ws.Range("E2")
$c.Value = "'  +3.73%  "
$c.Style = "Normal"

$c = $ws.Range("D3")
$c.Value = "'1.605.56"
$c.Style = "Normal"

$c = $ws.Range("E4")
$c.Value = "'  -0.07%  "
$c.Style = "Normal"

$c = $ws.Range("D5")
$c.Value = "'212.89"
$c.Style = "Normal"
$c = $ws.Range("E5")
$c.Value = "'  +1.00%  "
$c.Style = "Normal"

$c = $ws.Range("D6")
$c.Value = "'0.516"
$c.Style = "Normal"
$c = $ws.Range("E6")
$c.Value = "'  +6.77%  "
$c.Style = "Normal"

$c = $ws.Range("D7")
$c.Value = "'0.998"
$c.Style = "Normal"
$c = $ws.Range("E7")
$c.Value = "'  -0.20%  "
$c.Style = "Normal"

$c = $ws.Range("D8")
$c.Value = "'26.72"
$c.Style = "Normal"
$c = $ws.Range("E8")
$c.Value = "'  +10.56%  "
$c.Style = "Normal"

$c = $ws.Range("D9")
$c.Value = "'43.68"
$c.Style = "Normal"
$c = $ws.Range("E9")
$c.Value = "'  -0.66%  "
$c.Style = "Normal"

$c = $ws.Range("E10")
$c.Value = "'  +2.55%  "
$c.Style = "Normal"

$c = $ws.Range("D11")
$c.Value = "'0.0597"
$c.Style = "Normal"
$c = $ws.Range("E11")
$c.Value = "'  +2.38%  "
$c.Style = "Normal"

$c = $ws.Range("E12")
$c.Value = "'  +2.18%  "
$c.Style = "Normal"

$c = $ws.Range("D13")
$c.Value = "'1.835.56"
$c.Style = "Normal"
$c = $ws.Range("E13")
$c.Value = "'  +3.43%  "
$c.Style = "Normal"

$c = $ws.Range("D14")
$c.Value = "'1.627.55"
$c.Style = "Normal"
$c = $ws.Range("E14")
$c.Value = "'  +5.33%  "
$c.Style = "Normal"

$c = $ws.Range("D15")
$c.Value = "'29.550.11"
$c.Style = "Normal"
$c = $ws.Range("E15")
$c.Value = "'  +3.88%  "
$c.Style = "Normal"

$c = $ws.Range("E16")
$c.Value = "'  +5.44%  "
$c.Style = "Normal"

$c = $ws.Range("E17")
$c.Value = "'  +3.51%  "
$c.Style = "Normal"

$c = $ws.Range("D18")
$c.Value = "'63.42"
$c.Style = "Normal"
$c = $ws.Range("E18")
$c.Value = "'  +3.74%  "
$c.Style = "Normal"

$c = $ws.Range("D19")
$c.Value = "'240.15"
$c.Style = "Normal"
$c = $ws.Range("E19")
$c.Value = "'  +4.89%  "
$c.Style = "Normal"

$c = $ws.Range("D20")
$c.Value = "'7.58"
$c.Style = "Normal"
$c = $ws.Range("E20")
$c.Value = "'  +3.09%  "
$c.Style = "Normal"

$c = $ws.Range("E21")
$c.Value = "'  +2.99%  "
$c.Style = "Normal"

$c = $ws.Range("E22")
$c.Value = "'  -0.11%  "
$c.Style = "Normal"

$c = $ws.Range("E23")
$c.Value = "'  +3.47%  "
$c.Style = "Normal"

$c = $ws.Range("D24")
$c.Value = "'9.19"
$c.Style = "Normal"

$c = $ws.Range("E25")
$c.Value = "'  +0.25%  "
$c.Style = "Normal"

$c = $ws.Range("D26")
$c.Value = "'154.83"
$c.Style = "Normal"
$c = $ws.Range("E26")
$c.Value = "'  +2.56%  "
$c.Style = "Normal"

$c = $ws.Range("D27")
$c.Value = "'15.26"
$c.Style = "Normal"
$c = $ws.Range("E27")
$c.Value = "'  +3.60%  "
$c.Style = "Normal"

$c = $ws.Range("E28")
$c.Value = "'  +4.74%  "
$c.Style = "Normal"

$c = $ws.Range("D29")
$c.Value = "'6.36"
$c.Style = "Normal"
$c = $ws.Range("E29")
$c.Value = "'  +2.04%  "
$c.Style = "Normal"

$c = $ws.Range("D30")
$c.Value = "'0.999"
$c.Style = "Normal"
$c = $ws.Range("E30")
$c.Value = "'  -0.09%  "
$c.Style = "Normal"

$c = $ws.Range("D31")
$c.Value = "'0.0471"
$c.Style = "Normal"
$c = $ws.Range("E31")
$c.Value = "'  +1.08%  "
$c.Style = "Normal"

$c = $ws.Range("E32")
$c.Value = "'  +0.98%  "
$c.Style = "Normal"

$c = $ws.Range("E33")
$c.Value = "'  +2.59%  "
$c.Style = "Normal"

$c = $ws.Range("D34")
$c.Value = "'1.430.80"
$c.Style = "Normal"
$c = $ws.Range("E34")
$c.Value = "'  +3.03%  "
$c.Style = "Normal"

$c = $ws.Range("E35")
$c.Value = "'  +3.39%  "
$c.Style = "Normal"

$c = $ws.Range("E36")
$c.Value = "'  +0.74%  "
$c.Style = "Normal"

$c = $ws.Range("D37")
$c.Value = "'2.83"
$c.Style = "Normal"
$c = $ws.Range("E37")
$c.Value = "'  +6.65%  "
$c.Style = "Normal"

$c = $ws.Range("E38")
$c.Value = "'  +1.92%  "
$c.Style = "Normal"

$c = $ws.Range("E39")
$c.Value = "'  -0.36%  "
$c.Style = "Normal"

$c = $ws.Range("E40")
$c.Value = "'  +2.12%  "
$c.Style = "Normal"

$c = $ws.Range("D41")
$c.Value = "'0.532"
$c.Style = "Normal"
$c = $ws.Range("E41")
$c.Value = "'  +3.29%  "
$c.Style = "Normal"

$c = $ws.Range("D42")
$c.Value = "'1.94"
$c.Style = "Normal"
$c = $ws.Range("E42")
$c.Value = "'  -0.67%  "
$c.Style = "Normal"

$ws.Range("B43").Value = "BitcoinSV"
$ws.Range("C43").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$c = $ws.Range("D43")
$c.Value = "'53.13"
$c.Style = "Normal"
$c = $ws.Range("E43")
$c.Value = "'  +22.87%  "
$c.Style = "Normal"

$ws.Range("B44").Value = "PaxDollar"
$ws.Range("C44").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$c = $ws.Range("D44")
$c.Value = "'0.998"
$c.Style = "Normal"
$c = $ws.Range("E44")
$c.Value = "'  -0.15%  "
$c.Style = "Normal"

$ws.Range("B45").Value = "ARBITRUM"
$ws.Range("C45").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$c = $ws.Range("D45")
$c.Value = "'0.794"
$c.Style = "Normal"
$c = $ws.Range("E45")
$c.Value = "'  +2.91%  "
$c.Style = "Normal"

$c = $ws.Range("E46")
$c.Value = "'  +0.29%  "
$c.Style = "Normal"

$c = $ws.Range("D47")
$c.Value = "'65.83"
$c.Style = "Normal"
$c = $ws.Range("E47")
$c.Value = "'  +6.24%  "
$c.Style = "Normal"

$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$c = $ws.Range("D48")
$c.Value = "'5.27"
$c.Style = "Normal"
$c = $ws.Range("E48")
$c.Value = "'  -1.47%  "
$c.Style = "Normal"

$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$c = $ws.Range("D49")
$c.Value = "'1.745.72"
$c.Style = "Normal"
$c = $ws.Range("E49")
$c.Value = "'  +3.41%  "
$c.Style = "Normal"

$c = $ws.Range("D50")
$c.Value = "'86.84"
$c.Style = "Normal"
$c = $ws.Range("E50")
$c.Value = "'  +1.69%  "
$c.Style = "Normal"

$c = $ws.Range("D51")
$c.Value = "'0.838"
$c.Style = "Normal"
$c = $ws.Range("E51")
$c.Value = "'  -3.70%  "
$c.Style = "Normal"
